$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Julio de 2020 a las 22:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4539329
$ws.Range("C4").Value = 40986
$ws.Range("D4").Value = 2217963
$ws.Range("E4").Value = 2168212
$ws.Range("G4").Value = 863
$ws.Range("H4").Value = 153154

# Row 21 - Alemania
$ws.Range("B21").Value = 208811
$ws.Range("C21").Value = 860
$ws.Range("E21").Value = 7599
$ws.Range("G21").Value = 5
$ws.Range("H21").Value = 9212

# Row 64 - Uzbekistan
$ws.Range("B64").Value = 22585
$ws.Range("C64").Value = 692
$ws.Range("D64").Value = 12937
$ws.Range("E64").Value = 9517
$ws.Range("G64").Value = 5
$ws.Range("H64").Value = 131

# Rows 70-71 - Venezuela / Costa Rica swap order (Costa Rica now ranks above Venezuela)
# with Costa Rica receiving freshly updated figures, Venezuela keeping its prior totals.
$ws.Range("A70").Value = "Costa Rica"
$ws.Range("B70").Value = 16800
$ws.Range("C70").Value = 456
$ws.Range("D70").Value = 4050
$ws.Range("E70").Value = 12617
$ws.Range("G70").Value = 8
$ws.Range("H70").Value = 133

$ws.Range("A71").Value = "Venezuela"
$ws.Range("B71").Value = 16571
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 10195
$ws.Range("E71").Value = 6225
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 151

# Row 148 - Angola
$ws.Range("B148").Value = 1078
$ws.Range("C148").Value = 78
$ws.Range("D148").Value = 301
$ws.Range("E148").Value = 729
$ws.Range("G148").Value = 1
$ws.Range("H148").Value = 48

# Row 154 - Botsuana
$ws.Range("B154").Value = 804
$ws.Range("C154").Value = 65
$ws.Range("E154").Value = 739

# Rows 185-186 - Seychelles / San Martin (Parte Holandesa) swap order
# (San Martin now ranks above Seychelles) with San Martin receiving updated figures,
# Seychelles keeping its prior totals.
$ws.Range("A185").Value = "San Martin (Parte Holandesa)"
$ws.Range("B185").Value = 115
$ws.Range("C185").Value = 1
$ws.Range("D185").Value = 63
$ws.Range("E185").Value = 37
$ws.Range("H185").Value = 15

$ws.Range("A186").Value = "Seychelles"
$ws.Range("D186").Value = 39
$ws.Range("E186").Value = 75
$ws.Range("H186").Value = 0
